$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Havades_App" block (rows 20-23), mirroring the existing block at rows 16-19 ---

# Copy formatting (fill/border/font/alignment) from the existing block so the new
# rows reuse the same style definitions instead of creating new ones.
$ws.Range("A16:C19").Copy()
$ws.Range("A20").PasteSpecial(-4122)

# Category header row (merged A20:C20)
$ws.Range("A20").Value = "آذربایجان غربی ـ تغییر در روش ایمپورت اطلاعات بارگیری و فایل اکسل مورد نظر"
$ws.Range("A20:C20").Merge()

# "Havades_App" sub-block (merged A21:A23)
$ws.Range("A21").Value = "Havades_App"
$ws.Range("B21").Value = "frmBaseTableNotstd"
$ws.Range("B22").Value = "frmMPPostTransLoad"
$ws.Range("A21:A23").Merge()

# --- New "3rd Week" header row (row 25), mirroring rows 2 and 15 ---
$ws.Range("A2:C2").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value = "3rd Week (1400/2/4)"
$ws.Range("A25:C25").Merge()

$excel.CutCopyMode = 0

# --- Update sheet view (scrolled position / selection) ---
$ws.Range("C23").Select()
$ws.Application.ActiveWindow.ScrollRow = 10
